$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 328
$ws.Range("I2").Value = 328
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 328
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -215
$ws.Range("N2").Value = $null

$ws.Range("H7").Value = 12600
$ws.Range("J7").Value = 12600
$ws.Range("L7").Value = 12600
$ws.Range("N7").Value = -12824

$ws.Range("H14").Value = 12600
$ws.Range("J14").Value = 12600
$ws.Range("L14").Value = 12600
$ws.Range("N14").Value = -12982

$ws.Range("H76").Value = 1000
$ws.Range("I76").Value = 1000
$ws.Range("K76").Value = 1000
$ws.Range("M76").Value = -685

$ws.Range("H79").Value = 1000
$ws.Range("I79").Value = 1000
$ws.Range("K79").Value = 1000
$ws.Range("M79").Value = 92

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

$ws.Range("H113").Value = 2690
$ws.Range("I113").Value = 2251.6667
$ws.Range("J113").Value = 3347.5
$ws.Range("K113").Value = 2251.6667
$ws.Range("L113").Value = 3347.5
$ws.Range("M113").Value = 1002.3333
$ws.Range("N113").Value = -9855.5

$ws.Range("H116").Value = 2664.3333
$ws.Range("I116").Value = 1997
$ws.Range("K116").Value = 1997
$ws.Range("M116").Value = 1445

$ws.Range("H132").Value = 4648.5
$ws.Range("I132").Value = 1826.5714
$ws.Range("J132").Value = 8599.200000000001
$ws.Range("K132").Value = 5479.7142
$ws.Range("L132").Value = 25797.6
$ws.Range("M132").Value = -2949.7142
$ws.Range("N132").Value = -30857.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 25769.23
$ws.Range("I37").Value = 23545.455
$ws.Range("K37").Value = 23545.455
$ws.Range("M37").Value = -23272.455

$ws.Range("H45").Value = 5008.5557
$ws.Range("I45").Value = 3817.4
$ws.Range("K45").Value = 3817.4
$ws.Range("M45").Value = -3440.4

$ws.Range("H55").Value = 49500
$ws.Range("J55").Value = 49500
$ws.Range("L55").Value = 49500
$ws.Range("N55").Value = -50130

$ws.Range("H94").Value = 45999.332
$ws.Range("J94").Value = 45999.332
$ws.Range("L94").Value = 45999.332
$ws.Range("N94").Value = -47801.332

$ws.Range("H108").Value = 67852.28999999999
$ws.Range("J108").Value = 67852.28999999999
$ws.Range("L108").Value = 67852.28999999999
$ws.Range("N108").Value = -75532.28999999999

$ws.Range("H119").Value = 73975
$ws.Range("J119").Value = 73975
$ws.Range("L119").Value = 73975
$ws.Range("N119").Value = -83651

$ws.Range("H122").Value = 1892.5
$ws.Range("I122").Value = 1880.1666
$ws.Range("K122").Value = 5640.4998
$ws.Range("M122").Value = -3190.4998

$ws.Range("H132").Value = 1702.4546
$ws.Range("I132").Value = 1491.8889
$ws.Range("J132").Value = 2650
$ws.Range("K132").Value = 4475.6667
$ws.Range("L132").Value = 7950
$ws.Range("M132").Value = -1945.6667
$ws.Range("N132").Value = -13010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2201
$ws.Range("I86").Value = 2201
$ws.Range("K86").Value = 2201
$ws.Range("M86").Value = -1078

$ws.Range("H89").Value = 2201
$ws.Range("I89").Value = 2201
$ws.Range("K89").Value = 11005
$ws.Range("M89").Value = -5389

$ws.Range("H134").Value = 3743
$ws.Range("I134").Value = 1067
$ws.Range("K134").Value = 3201
$ws.Range("M134").Value = -666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 932.4286
$ws.Range("I16").Value = 921.1667
$ws.Range("K16").Value = 921.1667
$ws.Range("M16").Value = -634.1667

$ws.Range("H29").Value = 2500
$ws.Range("J29").Value = 2500
$ws.Range("L29").Value = 2500
$ws.Range("N29").Value = -3086

$ws.Range("H41").Value = 16079.8
$ws.Range("I41").Value = 3466.6667
$ws.Range("J41").Value = 34999.5
$ws.Range("K41").Value = 3466.6667
$ws.Range("L41").Value = 34999.5
$ws.Range("M41").Value = -3038.6667
$ws.Range("N41").Value = -35855.5

$ws.Range("H50").Value = 26666.334
$ws.Range("J50").Value = 34999.5
$ws.Range("L50").Value = 34999.5
$ws.Range("N50").Value = -36249.5

$ws.Range("H113").Value = 932.4286
$ws.Range("I113").Value = 921.1667
$ws.Range("K113").Value = 921.1667
$ws.Range("M113").Value = 1248.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 783
$ws.Range("I9").Value = 349
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 1047
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = -823
$ws.Range("N9").Value = -3448

$ws.Range("H80").Value = 7326.2
$ws.Range("J80").Value = 9394.571
$ws.Range("L80").Value = 28183.713
$ws.Range("N80").Value = -30055.713

$ws.Range("H83").Value = 7326.2
$ws.Range("J83").Value = 9394.571
$ws.Range("L83").Value = 84551.139
$ws.Range("N83").Value = -93911.139

$ws.Range("H123").Value = 5999.5
$ws.Range("J123").Value = 1999
$ws.Range("L123").Value = 5997
$ws.Range("N123").Value = -10897

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 15169.546
$ws.Range("J15").Value = 15169.546
$ws.Range("L15").Value = 15169.546
$ws.Range("N15").Value = -15745.546

$ws.Range("H43").Value = 3371.75
$ws.Range("I43").Value = 2326
$ws.Range("J43").Value = 3999.2
$ws.Range("K43").Value = 2326
$ws.Range("L43").Value = 3999.2
$ws.Range("M43").Value = -2175
$ws.Range("N43").Value = -4301.2

$ws.Range("H81").Value = 15169.546
$ws.Range("J81").Value = 15169.546
$ws.Range("L81").Value = 15169.546
$ws.Range("N81").Value = -17165.546

$ws.Range("H84").Value = 15169.546
$ws.Range("J84").Value = 15169.546
$ws.Range("L84").Value = 45508.638
$ws.Range("N84").Value = -55492.638

$ws.Range("H92").Value = 12876.429
$ws.Range("J92").Value = 12876.429
$ws.Range("L92").Value = 12876.429
$ws.Range("N92").Value = -16620.429

$ws.Range("H95").Value = 26562.666
$ws.Range("J95").Value = 26562.666
$ws.Range("L95").Value = 26562.666
$ws.Range("N95").Value = -32054.666

$ws.Range("H98").Value = 6000
$ws.Range("J98").Value = 6000
$ws.Range("L98").Value = 6000
$ws.Range("N98").Value = -11990

$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2109.2856
$ws.Range("I22").Value = 2041.1538
$ws.Range("K22").Value = 2041.1538
$ws.Range("M22").Value = -1746.1538

$ws.Range("H27").Value = 2109.2856
$ws.Range("I27").Value = 2041.1538
$ws.Range("K27").Value = 2041.1538
$ws.Range("M27").Value = -1934.1538

$ws.Range("H116").Value = 247250
$ws.Range("J116").Value = 247250
$ws.Range("L116").Value = 247250
$ws.Range("N116").Value = -256428

$ws.Range("H122").Value = 7004
$ws.Range("I122").Value = 7004
$ws.Range("K122").Value = 21012
$ws.Range("M122").Value = -18562

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 19992.5
$ws.Range("J31").Value = 19992.5
$ws.Range("L31").Value = 19992.5
$ws.Range("N31").Value = -20688.5

$ws.Range("H54").Value = 6038.5
$ws.Range("I54").Value = 5000
$ws.Range("J54").Value = 7077
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 7077
$ws.Range("M54").Value = -4480
$ws.Range("N54").Value = -8117

$ws.Range("H58").Value = 54999
$ws.Range("I58").Value = 54999
$ws.Range("K58").Value = 54999
$ws.Range("M58").Value = -54691

$ws.Range("H94").Value = 49999
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 49999
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 49999
$ws.Range("M94").Value = $null
$ws.Range("N94").Value = -51801

$ws.Range("H119").Value = 74832
$ws.Range("J119").Value = 74832
$ws.Range("L119").Value = 74832
$ws.Range("N119").Value = -84508

$ws.Range("H126").Value = 1257.3334
$ws.Range("J126").Value = 1124.5
$ws.Range("L126").Value = 3373.5
$ws.Range("N126").Value = -8313.5

$ws.Range("H132").Value = 250567.25
